# WebForm User Assignment execution
# Reassign the generated phone numbers (PN_Value, column F) and the
# Match1UserPos / Match2UserPos values (AM2 / AN2) on row 2 to the new
# run's values, per the latest test execution.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force the cell to stay text (shared-string) typed even though the
    # new value is all-digits — a plain `.Value = "123"` assignment would
    # otherwise be auto-coerced to a number by Excel.
    $range.NumberFormat = "@"
    $range.Value = $text
}

Set-TextValue $ws.Range("F2")  "9840014558"
Set-TextValue $ws.Range("F3")  "9840071961"
Set-TextValue $ws.Range("F4")  "9840011831"
Set-TextValue $ws.Range("F5")  "9840069917"
Set-TextValue $ws.Range("F6")  "9840027538"
Set-TextValue $ws.Range("F7")  "9840016153"
Set-TextValue $ws.Range("F8")  "9840070145"
Set-TextValue $ws.Range("F9")  "9840030416"
Set-TextValue $ws.Range("F10") "9840069697"

Set-TextValue $ws.Range("AM2") "1"
Set-TextValue $ws.Range("AN2") "0"
